$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2219423333333333
$ws.Range("H2").Value = 0.6658269999999999
$ws.Range("I2").Value = 0.01170735863810222
$ws.Range("J2").Value = 0.01170735863810222
$ws.Range("M2").Value = 176.8550973333333
$ws.Range("N2").Value = 530.565292
$ws.Range("O2").Value = 0.2669710696905332
$ws.Range("P2").Value = 0.2669710696905332
$ws.Range("Q2").Value = 39.25163296405378
$ws.Range("R2").Value = 353.264696676484
$ws.Range("S2").Value = 0.003125526058864854
$ws.Range("T2").Value = 0.003125526058864854
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2219423333333333
$ws.Range("H3").Value = 0.6658269999999999
$ws.Range("I3").Value = 0.01170735863810222
$ws.Range("J3").Value = 0.01170735863810222
$ws.Range("O3").Value = 0.2198524722701247
$ws.Range("P3").Value = 0.2198524722701247
$ws.Range("Q3").Value = 32.32398386008622
$ws.Range("R3").Value = 290.915854740776
$ws.Range("S3").Value = 0.002573891740339774
$ws.Range("T3").Value = 0.002573891740339774
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2219423333333333
$ws.Range("H4").Value = 0.6658269999999999
$ws.Range("I4").Value = 0.01170735863810222
$ws.Range("J4").Value = 0.01170735863810222
$ws.Range("M4").Value = 84.02511333333334
$ws.Range("N4").Value = 252.07534
$ws.Range("O4").Value = 0.1268398520919549
$ws.Range("P4").Value = 0.1268398520919549
$ws.Range("Q4").Value = 18.64872971179778
$ws.Range("R4").Value = 167.83856740618
$ws.Range("S4").Value = 0.001484959638044356
$ws.Range("T4").Value = 0.001484959638044356
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2219423333333333
$ws.Range("H5").Value = 0.6658269999999999
$ws.Range("I5").Value = 0.01170735863810222
$ws.Range("J5").Value = 0.01170735863810222
$ws.Range("M5").Value = 84.92877566666668
$ws.Range("N5").Value = 254.786327
$ws.Range("O5").Value = 0.1282039727953256
$ws.Range("P5").Value = 0.1282039727953256
$ws.Range("Q5").Value = 18.84929063860323
$ws.Range("R5").Value = 169.643615747429
$ws.Range("S5").Value = 0.001500929888344377
$ws.Range("T5").Value = 0.001500929888344377
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2219423333333333
$ws.Range("H6").Value = 0.6658269999999999
$ws.Range("I6").Value = 0.01170735863810222
$ws.Range("J6").Value = 0.01170735863810222
$ws.Range("M6").Value = 63.97102366666667
$ws.Range("N6").Value = 191.913071
$ws.Range("O6").Value = 0.09656726254996952
$ws.Range("P6").Value = 0.09656726254996952
$ws.Range("Q6").Value = 14.19787825830189
$ws.Range("R6").Value = 127.780904324717
$ws.Range("S6").Value = 0.001130547575372271
$ws.Range("T6").Value = 0.001130547575372271
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2219423333333333
$ws.Range("H7").Value = 0.6658269999999999
$ws.Range("I7").Value = 0.01170735863810222
$ws.Range("J7").Value = 0.01170735863810222
$ws.Range("M7").Value = 107.0290476666667
$ws.Range("N7").Value = 321.087143
$ws.Range("O7").Value = 0.1615653706020921
$ws.Range("P7").Value = 0.1615653706020921
$ws.Range("Q7").Value = 23.75427657358455
$ws.Range("R7").Value = 213.788489162261
$ws.Range("S7").Value = 0.001891503737136589
$ws.Range("T7").Value = 0.001891503737136589
$ws.Range("I8").Value = 0.8920180274758701
$ws.Range("J8").Value = 0.8920180274758701
$ws.Range("M8").Value = 176.8550973333333
$ws.Range("N8").Value = 530.565292
$ws.Range("O8").Value = 0.2669710696905332
$ws.Range("P8").Value = 0.2669710696905332
$ws.Range("Q8").Value = 2990.697158439299
$ws.Range("R8").Value = 26916.27442595369
$ws.Range("S8").Value = 0.2381430069784725
$ws.Range("T8").Value = 0.2381430069784725
$ws.Range("I9").Value = 0.8920180274758701
$ws.Range("J9").Value = 0.8920180274758701
$ws.Range("O9").Value = 0.2198524722701247
$ws.Range("P9").Value = 0.2198524722701247
$ws.Range("S9").Value = 0.1961123686500901
$ws.Range("T9").Value = 0.1961123686500901
$ws.Range("I10").Value = 0.8920180274758701
$ws.Range("J10").Value = 0.8920180274758701
$ws.Range("M10").Value = 84.02511333333334
$ws.Range("N10").Value = 252.07534
$ws.Range("O10").Value = 0.1268398520919549
$ws.Range("P10").Value = 0.1268398520919549
$ws.Range("Q10").Value = 1420.901469466307
$ws.Range("R10").Value = 12788.11322519676
$ws.Range("S10").Value = 0.1131434346683967
$ws.Range("T10").Value = 0.1131434346683967
$ws.Range("I11").Value = 0.8920180274758701
$ws.Range("J11").Value = 0.8920180274758701
$ws.Range("M11").Value = 84.92877566666668
$ws.Range("N11").Value = 254.786327
$ws.Range("O11").Value = 0.1282039727953256
$ws.Range("P11").Value = 0.1282039727953256
$ws.Range("Q11").Value = 1436.182795327076
$ws.Range("R11").Value = 12925.64515794368
$ws.Range("S11").Value = 0.1143602549274565
$ws.Range("T11").Value = 0.1143602549274565
$ws.Range("I12").Value = 0.8920180274758701
$ws.Range("J12").Value = 0.8920180274758701
$ws.Range("M12").Value = 63.97102366666667
$ws.Range("N12").Value = 191.913071
$ws.Range("O12").Value = 0.09656726254996952
$ws.Range("P12").Value = 0.09656726254996952
$ws.Range("Q12").Value = 1081.7780295117
$ws.Range("R12").Value = 9736.002265605295
$ws.Range("S12").Value = 0.08613973905856827
$ws.Range("T12").Value = 0.08613973905856827
$ws.Range("I13").Value = 0.8920180274758701
$ws.Range("J13").Value = 0.8920180274758701
$ws.Range("M13").Value = 107.0290476666667
$ws.Range("N13").Value = 321.087143
$ws.Range("O13").Value = 0.1615653706020921
$ws.Range("P13").Value = 0.1615653706020921
$ws.Range("Q13").Value = 1809.908074766212
$ws.Range("R13").Value = 16289.17267289591
$ws.Range("S13").Value = 0.1441192231928861
$ws.Range("T13").Value = 0.1441192231928861
$ws.Range("G14").Value = 1.825126666666667
$ws.Range("H14").Value = 5.475379999999999
$ws.Range("I14").Value = 0.09627461388602765
$ws.Range("J14").Value = 0.09627461388602765
$ws.Range("M14").Value = 176.8550973333333
$ws.Range("N14").Value = 530.565292
$ws.Range("O14").Value = 0.2669710696905332
$ws.Range("P14").Value = 0.2669710696905332
$ws.Range("Q14").Value = 322.7829542789955
$ws.Range("R14").Value = 2905.04658851096
$ws.Range("S14").Value = 0.02570253665319586
$ws.Range("T14").Value = 0.02570253665319586
$ws.Range("G15").Value = 1.825126666666667
$ws.Range("H15").Value = 5.475379999999999
$ws.Range("I15").Value = 0.09627461388602765
$ws.Range("J15").Value = 0.09627461388602765
$ws.Range("O15").Value = 0.2198524722701247
$ws.Range("P15").Value = 0.2198524722701247
$ws.Range("Q15").Value = 265.8139347726044
$ws.Range("R15").Value = 2392.32541295344
$ws.Range("S15").Value = 0.02116621187969486
$ws.Range("T15").Value = 0.02116621187969486
$ws.Range("G16").Value = 1.825126666666667
$ws.Range("H16").Value = 5.475379999999999
$ws.Range("I16").Value = 0.09627461388602765
$ws.Range("J16").Value = 0.09627461388602765
$ws.Range("M16").Value = 84.02511333333334
$ws.Range("N16").Value = 252.07534
$ws.Range("O16").Value = 0.1268398520919549
$ws.Range("P16").Value = 0.1268398520919549
$ws.Range("Q16").Value = 153.3564750143556
$ws.Range("R16").Value = 1380.2082751292
$ws.Range("S16").Value = 0.01221145778551381
$ws.Range("T16").Value = 0.01221145778551381
$ws.Range("G17").Value = 1.825126666666667
$ws.Range("H17").Value = 5.475379999999999
$ws.Range("I17").Value = 0.09627461388602765
$ws.Range("J17").Value = 0.09627461388602765
$ws.Range("M17").Value = 84.92877566666668
$ws.Range("N17").Value = 254.786327
$ws.Range("O17").Value = 0.1282039727953256
$ws.Range("P17").Value = 0.1282039727953256
$ws.Range("Q17").Value = 155.0057732365845
$ws.Range("R17").Value = 1395.05195912926
$ws.Range("S17").Value = 0.01234278797952477
$ws.Range("T17").Value = 0.01234278797952477
$ws.Range("G18").Value = 1.825126666666667
$ws.Range("H18").Value = 5.475379999999999
$ws.Range("I18").Value = 0.09627461388602765
$ws.Range("J18").Value = 0.09627461388602765
$ws.Range("M18").Value = 63.97102366666667
$ws.Range("N18").Value = 191.913071
$ws.Range("O18").Value = 0.09656726254996952
$ws.Range("P18").Value = 0.09656726254996952
$ws.Range("Q18").Value = 116.7552211879978
$ws.Range("R18").Value = 1050.79699069198
$ws.Range("S18").Value = 0.009296975916028972
$ws.Range("T18").Value = 0.009296975916028972
$ws.Range("G19").Value = 1.825126666666667
$ws.Range("H19").Value = 5.475379999999999
$ws.Range("I19").Value = 0.09627461388602765
$ws.Range("J19").Value = 0.09627461388602765
$ws.Range("M19").Value = 107.0290476666667
$ws.Range("N19").Value = 321.087143
$ws.Range("O19").Value = 0.1615653706020921
$ws.Range("P19").Value = 0.1615653706020921
$ws.Range("Q19").Value = 195.3415690043711
$ws.Range("R19").Value = 1758.07412103934
$ws.Range("S19").Value = 0.01555464367206938
$ws.Range("T19").Value = 0.01555464367206938
